$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first two menu rows (dates in column A) were updated to the new date.
# These are stored as plain text, not real dates, so force text format while
# writing the value, then restore the General format so the text sticks
# without Excel re-parsing it as a date.
$ws.Range("A2:A3").NumberFormat = "@"
$ws.Range("A2").Value = "05.04.2025"
$ws.Range("A3").Value = "05.04.2025"
$ws.Range("A2:A3").NumberFormat = "General"

# Move the cell selection to match where the editor left off.
$ws.Range("C9").Select()
